$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second mora period (2502) and the second worker (YINA DANNELYS
# MARTINEZ BARRIOS) - only one worker / one period remains (DANNA MILENA DIAZ
# SOLANO, periodo 2503). Deleting the rows shifts everything below them up
# (the old "signature" rows 23/24 become 21/22) and the now-unused shared
# strings are dropped automatically on save.
$ws.Rows("17:18").Delete()

# Updated Salario Basico for the remaining worker/period row
$ws.Range("G16").Value = 1649395

# Updated summary totals: Valor Mora total, worker count, period count
$ws.Range("E11").Value = 58940
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

Write-Output "Done"
